# Applies the cryptos.xlsx data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.399.06"
$ws.Range("D3").Value = "1.573.89"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "'291.09"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").Value = "'0.3764"
$ws.Range("E7").Value = "  +2.67%  "
$ws.Range("D8").Value = "'49.91"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'0.3414"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "'0.07668"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "'21.31"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "'5.981"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "'6.923"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "1.575.55"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "'90.61"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").Value = "'0.06750"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "'16.76"
$ws.Range("E21").Value = "  +2.59%  "
$ws.Range("D22").Value = "'6.232"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'0.5297"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "22.409.60"
$ws.Range("D26").Value = "'2.418"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").Value = "'2.752"
$ws.Range("E27").Value = "  -7.41%  "
$ws.Range("D28").Value = "'20.27"
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("D29").Value = "'145.33"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "'5.065"
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").Value = "'126.19"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("D32").Value = "1.754.07"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "'6.222"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").Value = "'1.015"
$ws.Range("E34").Value = "  +4.23%  "
$ws.Range("D35").Value = "'2.020"
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("D36").Value = "'10.10"
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("D37").Value = "'0.08528"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "'0.02561"
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("D39").Value = "'0.2319"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").Value = "'1.331"
$ws.Range("E40").Value = "  +5.44%  "
$ws.Range("D41").Value = "'0.06529"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'11.64"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.6494"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("D45").Value = "'14.23"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "'0.6041"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").Value = "'3.792"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "'1.308"
$ws.Range("E49").Value = "  +10.64%  "
$ws.Range("D50").Value = "'2.097"
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("D51").Value = "'125.73"
$ws.Range("E51").Value = "  +3.44%  "
